# [Refactor] 변수명 통일 중 - name을 displayName/EnemyName으로 통일, player를 character로 통일
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header "character_name" -> "character_displayName" (B1)
$ws.Range("B1").Value = "character_displayName"

# Widen column B to better fit the new, longer header text
$ws.Columns.Item(2).ColumnWidth = 19.9

# Move the active selection to D10, matching the saved view state
$ws.Range("D10").Select()
